$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the "NA" value from C10 (script found nothing to report this time)
$ws.Range("C10").Value = "'"

# Append new row of results from the script run
$ws.Range("A11").Value = "'2025-02-25"
$ws.Range("B11").Value = "Rien ne nous concerne aujourd'hui !"
$ws.Range("C11").Value = "NA"
$ws.Range("D11").Value = 112
